$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Objetivos value text (row 10, B/C) which now shows the real objectives text ---
$ws.Range("B10").Value = "Propiciar aos alunos os conhecimentos básicos de corrosão, nos aspectos termodinâmicos e cinéticos, e descrever as principais formas de ataque e as técnicas de proteção contra a corrosão e a oxidação metálica."
$ws.Range("C10").Value = "Propiciar aos alunos os conhecimentos básicos de corrosão, nos aspectos termodinâmicos e cinéticos, e descrever as principais formas de ataque e as técnicas de proteção contra a corrosão e a oxidação metálica."

# --- Clear old rows 13:25 entirely (removing stray cell nodes) before rebuilding ---
$ws.Range("A13:C25").Clear()

# Template cells carrying the three column styles (A=bold, B=wrap, C=red+wrap)
$templateA = $ws.Range("A3")
$templateB = $ws.Range("B3")
$templateC = $ws.Range("C3")

# Row 13
$templateB.Copy($ws.Range("B13"))
$ws.Range("B13").Value = "5817344 - Livia Melo Carneiro"
$templateC.Copy($ws.Range("C13"))
$ws.Range("C13").Value = "5817344 - Livia Melo Carneiro"
$ws.Rows.Item(13).AutoFit()

# Row 14
$templateA.Copy($ws.Range("A14"))
$ws.Range("A14").Value = "Programa resumido:"
$templateB.Copy($ws.Range("B14"))
$ws.Range("B14").Value = "1. Princípios da corrosão. 2. Cinética da corrosão. 3. Formas de corrosão. 4. Proteção contra a corrosão. 5. Oxidação em temperaturas elevadas."
$templateC.Copy($ws.Range("C14"))
$ws.Range("C14").Value = "1. Princípios da corrosão. 2. Cinética da corrosão. 3. Formas de corrosão. 4. Proteção contra a corrosão. 5. Oxidação em temperaturas elevadas."
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$templateA.Copy($ws.Range("A15"))
$ws.Range("A15").Value = "Short syllabus:"
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$templateA.Copy($ws.Range("A16"))
$ws.Range("A16").Value = "Programa:"
$templateB.Copy($ws.Range("B16"))
$ws.Range("B16").Value = "1. Princípios da corrosão: Reações de oxi-redução. Potenciais de eletrodo - Sistema redox em estado de equilíbrio - Diagrama de Pourbaix`n2. Cinética da corrosão: - Sistema redox em estado de não equilíbrio - Teoria do potencial misto  Passivação.`n3. Formas de corrosão: - Corrosão galvânica - Corrosão por pites e frestas - Corrosão intergranular - Corrosão sob tensão - Danos causados pelo hidrogênio.`n4. Proteção contra a corrosão: - Proteção catódica e anódica - Inibidores  Revestimentos.`n5. Oxidação em temperaturas elevadas - Fundamentos termodinâmicos - Mecanismos de transporte - Velocidade de oxidação - Oxidação de metais puros - Oxidação de ligas."
$templateC.Copy($ws.Range("C16"))
$ws.Range("C16").Value = "1. Princípios da corrosão: Reações de oxi-redução. Potenciais de eletrodo - Sistema redox em estado de equilíbrio - Diagrama de Pourbaix`n2. Cinética da corrosão: - Sistema redox em estado de não equilíbrio - Teoria do potencial misto  Passivação.`n3. Formas de corrosão: - Corrosão galvânica - Corrosão por pites e frestas - Corrosão intergranular - Corrosão sob tensão - Danos causados pelo hidrogênio.`n4. Proteção contra a corrosão: - Proteção catódica e anódica - Inibidores  Revestimentos.`n5. Oxidação em temperaturas elevadas - Fundamentos termodinâmicos - Mecanismos de transporte - Velocidade de oxidação - Oxidação de metais puros - Oxidação de ligas."
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$templateA.Copy($ws.Range("A17"))
$ws.Range("A17").Value = "Syllabus:"
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$templateA.Copy($ws.Range("A18"))
$ws.Range("A18").Value = "Avaliação:"
$ws.Rows.Item(18).AutoFit()

# Row 19
$templateA.Copy($ws.Range("A19"))
$ws.Range("A19").Value = "Método:"
$templateB.Copy($ws.Range("B19"))
$ws.Range("B19").Value = "O aluno será avaliado através de duas provas escritas P1 e P2."
$templateC.Copy($ws.Range("C19"))
$ws.Range("C19").Value = "O aluno será avaliado através de duas provas escritas P1 e P2."
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$templateA.Copy($ws.Range("A20"))
$ws.Range("A20").Value = "Critério:"
$templateB.Copy($ws.Range("B20"))
$ws.Range("B20").Value = "A nota final NF será (P1 + P2)/2 ."
$templateC.Copy($ws.Range("C20"))
$ws.Range("C20").Value = "A nota final NF será (P1 + P2)/2 ."
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$templateA.Copy($ws.Range("A21"))
$ws.Range("A21").Value = "Norma de recuperação:"
$templateB.Copy($ws.Range("B21"))
$ws.Range("B21").Value = "Prova escrita sobre toda matéria. A média final MF será a média da nota final NF e da nota obtida na recuperação NR: MF = (NF + NR)/2 ."
$templateC.Copy($ws.Range("C21"))
$ws.Range("C21").Value = "Prova escrita sobre toda matéria. A média final MF será a média da nota final NF e da nota obtida na recuperação NR: MF = (NF + NR)/2 ."
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$templateA.Copy($ws.Range("A22"))
$ws.Range("A22").Value = "Bibliografia:"
$templateB.Copy($ws.Range("B22"))
$ws.Range("B22").Value = "V.GENTIL, Corrosão, Ed. Guanabara Dois, 1982`nL.V. RAMANATHAN, Corrosão e seu controle, Ed. Hermes`nL.L. SHREIR, Corrosion, Newnes Butterworths, 2 vol., 1976`nN. BIRKS and G.H.MEIER, Introduction to High Temperature Oxidation of Metals, Edward Arnold, 1983"
$templateC.Copy($ws.Range("C22"))
$ws.Range("C22").Value = "V.GENTIL, Corrosão, Ed. Guanabara Dois, 1982`nL.V. RAMANATHAN, Corrosão e seu controle, Ed. Hermes`nL.L. SHREIR, Corrosion, Newnes Butterworths, 2 vol., 1976`nN. BIRKS and G.H.MEIER, Introduction to High Temperature Oxidation of Metals, Edward Arnold, 1983"
$ws.Rows.Item(22).RowHeight = 120

# Row 23
$templateA.Copy($ws.Range("A23"))
$ws.Range("A23").Value = "Requisitos:"
$ws.Rows.Item(23).AutoFit()

# Row 24
$templateB.Copy($ws.Range("B24"))
$ws.Range("B24").Value = "LOB1053 -  Física III  (Requisito fraco)`n"
$templateC.Copy($ws.Range("C24"))
$ws.Range("C24").Value = "LOB1053 -  Física III  (Requisito fraco)`n"
$ws.Rows.Item(24).RowHeight = 30

# Row 25
$templateB.Copy($ws.Range("B25"))
$ws.Range("B25").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n"
$templateC.Copy($ws.Range("C25"))
$ws.Range("C25").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n"
$ws.Rows.Item(25).RowHeight = 30

# --- Adjust column definitions: narrow col A range from (1,2) to (1,1) ---
$ws.Range("B:B").ColumnWidth = 60.7109375
